# Update the Handback status report timestamps (re-generated report dates).
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for the first file row.
$wsOverview.Range("G2").Value = "2016-08-18 19:09:31"

# zh-cn sheet: "Correspond Handoff Datetime" (H) and "Correspond Handback DateTime" (K)
$wsZhCn.Range("H2").Value = "2016-08-18 19:09:26"
$wsZhCn.Range("K2").Value = "2016-08-18 19:09:53"

# de-de sheet: "Correspond Handback DateTime" (K)
$wsDeDe.Range("K2").Value = "2016-08-18 19:10:06"
